# Recreate the merges: rename header labels to reflect the new merge-suffix
# scheme, and swap the N/O (seas_id_y / season_ending_year_y) column data so
# that N holds the season-ending calendar year (as text) and O holds the new
# player_id_y values produced by the re-run merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row relabeling ------------------------------------------------
$ws.Range("K1").Value = "seas_id"
$ws.Range("L1").Value = "player_id_x"
$ws.Range("N1").Value = "season_ending_year_y"
$ws.Range("O1").Value = "player_id_y"

# ---- Per-row data swap (rows 2..70) ---------------------------------------
# New O values (player_id_y) from the re-run merge, one per data row.
$newPlayerIdY = @(3872,2711,3872,3872,1947,1947,2354,4450,4666,4666,3089,3260,3260,1344,3260,3260,3149,1431,4703,4703,3092,4816,4816,132,4575,2977,3643,2977,3643,1239,2035,803,3643,3643,3403,3403,3643,3403,3209,3209,3209,3785,3785,2932,2974,3785,422,2974,2974,512,2974,1198,2974,2974,5219,5162,5224,5224,5224,415,3934,415,415,415,5224,523,415,477,523)

for ($i = 0; $i -lt $newPlayerIdY.Length; $i++) {
    $row = $i + 2

    # N<row>: season_ending_year_y becomes the calendar year stored as TEXT
    # (same digits as the existing calendar_year/AY column), not a number.
    $yearText = [string]($ws.Cells.Item($row, 51).Value2)
    $ws.Range("ZZ1").Formula = "=""" + $yearText + """"
    $ws.Range("ZZ1").Copy()
    $ws.Cells.Item($row, 14).PasteSpecial(-4163)

    # O<row>: player_id_y becomes the new numeric id from the re-run merge.
    $ws.Cells.Item($row, 15).Value = $newPlayerIdY[$i]
}

# Clean up the scratch cell used to coerce text-typed literals.
$ws.Range("ZZ1").ClearContents()
